$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 4 for the incoming "Sponsoraanvraag" mail ---
$logs.Range("A4").Value = "Sponsoraanvraag"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Zou uw bedrijf bereid zijn om ons sportevenement te sponsoren?"
$logs.Range("D4").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F4").Value = "2025-06-20 09:00:10"
$logs.Range("G4").Value = "Nee"

# Extend the conditional formatting ranges on the Logs sheet to include row 4
$fcD = $logs.Range("D2:D3").FormatConditions
for ($i = 1; $i -le $fcD.Count; $i++) {
    $fcD.Item($i).ModifyAppliesToRange($logs.Range("D2:D4"))
}

$fcG = $logs.Range("G2:G3").FormatConditions
for ($i = 1; $i -le $fcG.Count; $i++) {
    $fcG.Item($i).ModifyAppliesToRange($logs.Range("G2:G4"))
}

# --- Dashboard sheet: new category tally row ---
$dash.Range("A4").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B4").Value = 1

# --- Update the chart series so it covers the new Dashboard row ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
